# adding remote execution feature using docker
# - RUNMANAGER: flip the leaveFeatureTest "execute" flag from "yes" to "Yes"
# - DATA: switch the browser under test to "firefox" for the login tests
#   (both the valid-creds and invalid-creds rows) to support remote/docker execution

$wb = $excel.ActiveWorkbook

$runManager = $wb.Worksheets.Item("RUNMANAGER")
$data = $wb.Worksheets.Item("DATA")

# RUNMANAGER: leaveFeatureTest row - execute flag "yes" -> "Yes"
$runManager.Range("C4").Value = "Yes"

# DATA sheet: switch browser column to firefox for the two login test rows
$data.Range("C2").Value = "firefox"
$data.Range("C3").Value = "firefox"

# Update selections to match - select RUNMANAGER's cell first, then DATA's,
# so DATA ends up as the active/selected tab (as in the original workbook).
$runManager.Range("C4").Select() | Out-Null
$data.Range("F15").Select() | Out-Null
